# Update odds values in row 3 (sheet "Sheet1") of the FlashScore workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value  = 1.6
$ws.Range("H3").Value  = 3.8
$ws.Range("I3").Value  = 5.75
$ws.Range("J3").Value  = 2.2
$ws.Range("K3").Value  = 2.25
$ws.Range("L3").Value  = 5.5
$ws.Range("M3").Value  = 1.05
$ws.Range("N3").Value  = 11
$ws.Range("O3").Value  = 1.25
$ws.Range("P3").Value  = 3.75
$ws.Range("Q3").Value  = 1.85
$ws.Range("U3").Value  = 1.83
$ws.Range("V3").Value  = 1.83
$ws.Range("W3").Value  = 7
$ws.Range("X3").Value  = 7.5
$ws.Range("Z3").Value  = 12
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 11
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 251
$ws.Range("AI3").Value = 29
$ws.Range("AJ3").Value = 17
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 8
$ws.Range("AR3").Value = 51
$ws.Range("AT3").Value = 3
$ws.Range("AU3").Value = 8.5
$ws.Range("AX3").Value = 7
$ws.Range("AY3").Value = 29
$ws.Range("AZ3").Value = 34
$ws.Range("BA3").Value = 101
$ws.Range("BB3").Value = 126
$ws.Range("BC3").Value = 251
